$wb = $excel.ActiveWorkbook

# --- Sheet 1: Significant Components ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = "['PPUNIT' 'QFHH' 'QNOHLTH' 'QSERV' 'QEXTRCT' 'QESL' 'QHISPC' 'QEDLESHI'`n 'PERCAP']"
$ws.Range("C3").Value = "['PERCAP' 'QRICH' 'MDHSEVAL']"
$ws.Range("C4").Value = "['PPUNIT' 'QRENTER' 'QFAM' 'QNOAUTO' 'QPOVTY']"
$ws.Range("C5").Value = "['QFEMALE' 'QAGEDEP' 'QFEMLBR']"

# --- Sheet 2: Loading Factors ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 0.7902275165323367
$ws.Range("C2").Value = -0.02455912210047059
$ws.Range("D2").Value = -0.4600798265337037
$ws.Range("E2").Value = 0.06838368148083136
$ws.Range("F2").Value = -0.1573190497643179
$ws.Range("A3").Value = "QFHH"
$ws.Range("B3").Value = 0.6130097814099174
$ws.Range("C3").Value = 0.2399781594692201
$ws.Range("D3").Value = 0.1790883810539986
$ws.Range("E3").Value = 0.262013816056053
$ws.Range("F3").Value = -0.09687620126079302
$ws.Range("A4").Value = "QNOHLTH"
$ws.Range("B4").Value = 0.6164840039903862
$ws.Range("C4").Value = 0.445338319500932
$ws.Range("D4").Value = 0.3253808145275445
$ws.Range("E4").Value = -0.2013657287527142
$ws.Range("F4").Value = -0.0821058418095862
$ws.Range("B5").Value = 0.5381204752784663
$ws.Range("C5").Value = 0.3662181480252166
$ws.Range("D5").Value = 0.29766950120097
$ws.Range("E5").Value = 0.001243120973004728
$ws.Range("F5").Value = -0.1694031916490242
$ws.Range("A6").Value = "QEXTRCT"
$ws.Range("B6").Value = 0.7251699864974676
$ws.Range("C6").Value = 0.1597478186990243
$ws.Range("D6").Value = 0.137735901008178
$ws.Range("E6").Value = -0.2661556141792868
$ws.Range("F6").Value = 0.0542731569191852
$ws.Range("A7").Value = "QESL"
$ws.Range("B7").Value = 0.8069195274027448
$ws.Range("C7").Value = 0.1431421462673804
$ws.Range("D7").Value = 0.2282344066406292
$ws.Range("E7").Value = -0.2351410662516274
$ws.Range("F7").Value = -0.0215068959427597
$ws.Range("A8").Value = "QHISPC"
$ws.Range("B8").Value = 0.8388032609392188
$ws.Range("C8").Value = 0.3352021757011634
$ws.Range("D8").Value = 0.1601474932357065
$ws.Range("E8").Value = -0.06922549950080999
$ws.Range("F8").Value = -0.09574193664954313
$ws.Range("A9").Value = "QEDLESHI"
$ws.Range("B9").Value = 0.8855973742225899
$ws.Range("C9").Value = 0.2225892957832175
$ws.Range("D9").Value = 0.2488822041126817
$ws.Range("E9").Value = -0.1306950260795137
$ws.Range("F9").Value = 0.01576001153787721
$ws.Range("B10").Value = 0.4725767454152621
$ws.Range("C10").Value = 0.7137770005873879
$ws.Range("D10").Value = 0.2336096424972315
$ws.Range("E10").Value = 0.005326542191881722
$ws.Range("F10").Value = -0.2730791484659
$ws.Range("A11").Value = "QRICH"
$ws.Range("B11").Value = 0.1589560307154843
$ws.Range("C11").Value = 0.858248130084476
$ws.Range("D11").Value = 0.3041837792558107
$ws.Range("E11").Value = 0.004009561501261934
$ws.Range("F11").Value = -0.1338394773551382
$ws.Range("A12").Value = "MDHSEVAL"
$ws.Range("B12").Value = 0.3726973490978535
$ws.Range("C12").Value = 0.8017743578204324
$ws.Range("D12").Value = 0.08260669060806446
$ws.Range("E12").Value = 0.03067618575093604
$ws.Range("F12").Value = -0.06162804297150613
$ws.Range("A13").Value = "QRENTER"
$ws.Range("B13").Value = -0.007034391187877889
$ws.Range("C13").Value = 0.2082000118516525
$ws.Range("D13").Value = 0.7800125732977286
$ws.Range("E13").Value = -0.1101264153554879
$ws.Range("F13").Value = -0.4193040231013465
$ws.Range("A14").Value = "QFAM"
$ws.Range("B14").Value = 0.2493636374503665
$ws.Range("C14").Value = 0.2498208255995001
$ws.Range("D14").Value = 0.5200879797823313
$ws.Range("E14").Value = 0.09552580902966797
$ws.Range("F14").Value = -0.1315017286932762
$ws.Range("A15").Value = "QNOAUTO"
$ws.Range("B15").Value = 0.1911459589363538
$ws.Range("C15").Value = 0.09747168808006407
$ws.Range("D15").Value = 0.6356287811780605
$ws.Range("E15").Value = -0.006055161056030798
$ws.Range("F15").Value = -0.03646340148529349
$ws.Range("A16").Value = "QPOVTY"
$ws.Range("B16").Value = 0.4090324310498006
$ws.Range("C16").Value = 0.1933683360348016
$ws.Range("D16").Value = 0.4940719745037385
$ws.Range("E16").Value = -0.02947174046105295
$ws.Range("F16").Value = -0.3383025193834173
$ws.Range("A17").Value = "QFEMALE"
$ws.Range("B17").Value = -0.06658500548861797
$ws.Range("C17").Value = -0.06807863240491877
$ws.Range("D17").Value = -0.0204800386068439
$ws.Range("E17").Value = 0.8696948642610516
$ws.Range("F17").Value = 0.1203783700291517
$ws.Range("A18").Value = "QAGEDEP"
$ws.Range("B18").Value = 0.01678028762196785
$ws.Range("C18").Value = -0.1306700629132541
$ws.Range("D18").Value = -0.09147423576591124
$ws.Range("E18").Value = 0.6945701925506927
$ws.Range("F18").Value = 0.5699008951874491
$ws.Range("A19").Value = "QFEMLBR"
$ws.Range("B19").Value = -0.1982536588688737
$ws.Range("C19").Value = 0.1418071559978814
$ws.Range("D19").Value = 0.0386822705246177
$ws.Range("E19").Value = 0.7433017059610204
$ws.Range("F19").Value = -0.02001253565263549
$ws.Range("B20").Value = -0.3244474135117502
$ws.Range("C20").Value = -0.2578790660093226
$ws.Range("D20").Value = -0.2749783774492807
$ws.Range("E20").Value = -0.001848682738725
$ws.Range("F20").Value = 0.8112807183109635
$ws.Range("B21").Value = 0.05321156903320884
$ws.Range("C21").Value = -0.04844388727422573
$ws.Range("D21").Value = -0.1337452409546111
$ws.Range("E21").Value = 0.1584510221529807
$ws.Range("F21").Value = 0.750291967498725

# --- Sheet 3: All Refactor Variances ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 4.524595077502688
$ws.Range("C2").Value = 3.113165441052899
$ws.Range("D2").Value = 2.536632550782927
$ws.Range("E2").Value = 2.093348638054015
$ws.Range("F2").Value = 2.093236337350453
$ws.Range("G2").Value = 1.983034029122986
$ws.Range("H2").Value = 0.8849929909863777
$ws.Range("I2").Value = 5.029294805889601
$ws.Range("J2").Value = 3.088512148077809
$ws.Range("K2").Value = 2.341179479818402
$ws.Range("L2").Value = 2.102363009327097
$ws.Range("M2").Value = 2.0651993393022
$ws.Range("N2").Value = 5.140423323011532
$ws.Range("O2").Value = 2.750776522408552
$ws.Range("P2").Value = 2.38162906535214
$ws.Range("Q2").Value = 2.101428973837291
$ws.Range("R2").Value = 2.048359020620988
$ws.Range("B3").Value = 0.1675775954630625
$ws.Range("C3").Value = 0.1153024237427
$ws.Range("D3").Value = 0.093949353732701
$ws.Range("E3").Value = 0.0775314310390376
$ws.Range("F3").Value = 0.07752727175372048
$ws.Range("G3").Value = 0.0734457047823328
$ws.Range("H3").Value = 0.03277751818468066
$ws.Range("I3").Value = 0.2394902288518857
$ws.Range("J3").Value = 0.1470720070513242
$ws.Range("K3").Value = 0.1114847371342096
$ws.Range("L3").Value = 0.1001125242536713
$ws.Range("M3").Value = 0.09834282568105714
$ws.Range("N3").Value = 0.2570211661505766
$ws.Range("O3").Value = 0.1375388261204276
$ws.Range("P3").Value = 0.119081453267607
$ws.Range("Q3").Value = 0.1050714486918645
$ws.Range("R3").Value = 0.1024179510310494
$ws.Range("B4").Value = 0.1675775954630625
$ws.Range("C4").Value = 0.2828800192057625
$ws.Range("D4").Value = 0.3768293729384635
$ws.Range("E4").Value = 0.4543608039775011
$ws.Range("F4").Value = 0.5318880757312217
$ws.Range("G4").Value = 0.6053337805135545
$ws.Range("H4").Value = 0.6381112986982351
$ws.Range("I4").Value = 0.2394902288518857
$ws.Range("J4").Value = 0.3865622359032099
$ws.Range("K4").Value = 0.4980469730374195
$ws.Range("L4").Value = 0.5981594972910907
$ws.Range("M4").Value = 0.6965023229721479
$ws.Range("N4").Value = 0.2570211661505766
$ws.Range("O4").Value = 0.3945599922710042
$ws.Range("P4").Value = 0.5136414455386112
$ws.Range("Q4").Value = 0.6187128942304758
$ws.Range("R4").Value = 0.7211308452615252
$ws.Range("B5").Value = 0.2626149949153471
$ws.Range("C5").Value = 0.1806932802755885
$ws.Range("D5").Value = 0.1472303560904819
$ws.Range("E5").Value = 0.121501423336656
$ws.Range("F5").Value = 0.1214949052177548
$ws.Range("G5").Value = 0.1150985806585216
$ws.Range("H5").Value = 0.05136645950565005
$ws.Range("I5").Value = 0.3438469922539836
$ws.Range("J5").Value = 0.2111579562631343
$ws.Range("K5").Value = 0.1600636975027975
$ws.Range("L5").Value = 0.1437360952745518
$ws.Range("M5").Value = 0.141195258705533
$ws.Range("N5").Value = 0.3564140514019552
$ws.Range("O5").Value = 0.1907265886963243
$ws.Range("P5").Value = 0.1651315486642663
$ws.Range("Q5").Value = 0.1457037226770675
$ws.Range("R5").Value = 0.1420240885603867

# --- Sheet 4: Final Variances ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 5.140423323011532
$ws.Range("C2").Value = 2.750776522408552
$ws.Range("D2").Value = 2.38162906535214
$ws.Range("E2").Value = 2.101428973837291
$ws.Range("F2").Value = 2.048359020620988
$ws.Range("B3").Value = 0.2570211661505766
$ws.Range("C3").Value = 0.1375388261204276
$ws.Range("D3").Value = 0.119081453267607
$ws.Range("E3").Value = 0.1050714486918645
$ws.Range("F3").Value = 0.1024179510310494
$ws.Range("B4").Value = 0.2570211661505766
$ws.Range("C4").Value = 0.3945599922710042
$ws.Range("D4").Value = 0.5136414455386112
$ws.Range("E4").Value = 0.6187128942304758
$ws.Range("F4").Value = 0.7211308452615252
$ws.Range("B5").Value = 0.3564140514019552
$ws.Range("C5").Value = 0.1907265886963243
$ws.Range("D5").Value = 0.1651315486642663
$ws.Range("E5").Value = 0.1457037226770675
$ws.Range("F5").Value = 0.1420240885603867

# --- Sheet 5: Included and Excluded ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = "[['PPUNIT', 'QFHH', 'QNOHLTH', 'QSERV', 'QEXTRCT', 'QESL', 'QHISPC', 'QEDLESHI', 'PERCAP', 'QRICH', 'MDHSEVAL', 'QRENTER', 'QFAM', 'QNOAUTO', 'QPOVTY', 'QFEMALE', 'QAGEDEP', 'QFEMLBR', 'MEDAGE', 'QSSBEN']]"

